# Penambahan dan Update API Engine, Dokumentasi API Engine, serta Penambahan Model Database
#
# Inserts a new catalogue row for the "Bill Of Material Detail" synchronize
# API, directly above the existing "Material Product Assembly" row (which
# currently sits at row 214). Every row from 214 downward shifts down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 214 (pushes 214..284 down to 215..285).
$ws.Rows(214).EntireRow.Insert()

# Clone formatting from the row directly below (the row that used to be 214,
# now 215) so the new row's cell styles match the rest of the table exactly
# instead of the engine synthesizing brand-new style entries.
$ws.Range("A215:H215").Copy()
$ws.Range("A214:H214").PasteSpecial(-4122)

# Populate the new row with the new API catalogue entry.
$ws.Range("B214").Value = "transaction.synchronize.production.setBillOfMaterialDetail"
$ws.Range("C214").Value = "Menyinkronkan Data Bill Of Material Detail"

# Reflect the author's final cursor position/selection in the sheet.
$ws.Range("C214").Select()
